# Update cryptos list sheet with latest values (prices / 1h volume %) and
# swap two pairs of rows (Toncoin/Monero and HuobiToken/RocketPoolETH)
# that changed rank order.
#
# Several "Price" values look like plain decimal numbers (e.g. "304.06"),
# but the source sheet always stores them as text. Force those specific
# cells to Text format first so Excel doesn't silently convert them to
# numbers when the new value is assigned.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textPriceCells = @("D5","D9","D10","D11","D12","D14","D17","D19","D22","D23","D27","D28","D29","D31","D33","D34","D39","D43","D45","D49","D50","D51")
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "43.177.91"
$ws.Range("E2").Value = "  +0.21%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.328.97"
$ws.Range("E3").Value = "  +0.88%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
$ws.Range("D5").Value = "304.06"
$ws.Range("E5").Value = "  +1.31%  "

# Row 6 - Solana
$ws.Range("E6").Value = "  -0.68%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -1.41%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.06%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.502"
$ws.Range("E9").Value = "  -1.23%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "35.57"
$ws.Range("E10").Value = "  -1.57%  "

# Row 11 - Chainlink
$ws.Range("D11").Value = "19.59"
$ws.Range("E11").Value = "  +7.98%  "

# Row 12 - Dogecoin
$ws.Range("D12").Value = "0.0799"
$ws.Range("E12").Value = "  +0.94%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +0.35%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "6.95"
$ws.Range("E14").Value = "  +2.04%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.692.43"
$ws.Range("E15").Value = "  +0.98%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.336.41"
$ws.Range("E16").Value = "  +1.00%  "

# Row 17 - Polygon
$ws.Range("D17").Value = "0.787"
$ws.Range("E17").Value = "  +0.48%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "43.005.11"
$ws.Range("E18").Value = "  -0.01%  "

# Row 19 - InternetComputer(DFINITY)
$ws.Range("D19").Value = "12.65"
$ws.Range("E19").Value = "  -1.17%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  -0.52%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +0.35%  "

# Row 22 - Litecoin
$ws.Range("D22").Value = "67.92"
$ws.Range("E22").Value = "  -0.22%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "237.08"
$ws.Range("E23").Value = "  -1.49%  "

# Row 24 - ImmutableX
$ws.Range("E24").Value = "  +3.03%  "

# Row 25 - Dai
$ws.Range("E25").Value = "  +0.12%  "

# Row 26 - PancakeSwap
$ws.Range("E26").Value = "  -0.18%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "24.95"
$ws.Range("E27").Value = "  -1.86%  "

# Row 28/29 swap - Toncoin and Monero traded ranking positions
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "165.70"
$ws.Range("E28").Value = "  -0.06%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.06"
$ws.Range("E29").Value = "  +2.04%  "

# Row 30 - Cosmos
$ws.Range("E30").Value = "  +0.60%  "

# Row 31 - InjectiveProtocol
$ws.Range("D31").Value = "33.18"
$ws.Range("E31").Value = "  -0.35%  "

# Row 32 - FirstDigitalUSD
$ws.Range("E32").Value = "  +0.00%  "

# Row 33 - Celestia
$ws.Range("D33").Value = "18.19"
$ws.Range("E33").Value = "  +6.31%  "

# Row 34 - Filecoin
$ws.Range("D34").Value = "5.00"
$ws.Range("E34").Value = "  -0.88%  "

# Row 35 - RenderToken
$ws.Range("E35").Value = "  -8.12%  "

# Row 36 - WEMIXToken
$ws.Range("E36").Value = "  -1.28%  "

# Row 37 - Hedera
$ws.Range("E37").Value = "  +1.29%  "

# Row 38 - Kaspa
$ws.Range("E38").Value = "  +0.33%  "

# Row 39 - LidoDAOToken
$ws.Range("D39").Value = "2.79"
$ws.Range("E39").Value = "  +1.32%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  -0.27%  "

# Row 42 - Maker
$ws.Range("D42").Value = "1.993.91"
$ws.Range("E42").Value = "  -1.33%  "

# Row 43 - FraxShare
$ws.Range("D43").Value = "10.73"
$ws.Range("E43").Value = "  +6.05%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  -0.52%  "

# Row 45 - EnergySwap
$ws.Range("D45").Value = "18.00"
$ws.Range("E45").Value = "  +2.60%  "

# Row 46 - ApeXProtocol
$ws.Range("E46").Value = "  -3.91%  "

# Row 47 - NEARProtocol
$ws.Range("E47").Value = "  -0.92%  "

# Row 48/49 swap - HuobiToken and RocketPoolETH traded ranking positions
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "2.559.98"
$ws.Range("E48").Value = "  +0.98%  "

$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").Value = "2.85"
$ws.Range("E49").Value = "  -1.93%  "

# Row 50 - MultiversX
$ws.Range("D50").Value = "53.61"
$ws.Range("E50").Value = "  -0.70%  "

# Row 51 - BitcoinSV
$ws.Range("D51").Value = "71.95"
$ws.Range("E51").Value = "  -0.87%  "
